$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Sheet1"

# Fill in the names
$ws.Range("A1").Value = "Jelena"
$ws.Range("B1").Value = "Vasilijevic"
$ws.Range("A2").Value = "Pera"
$ws.Range("B2").Value = "Peric"
